# Recalculated Leve profit figures (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets to
# reflect refreshed market-board pricing data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 126 - Saigaskin Codex
$ws.Range("H126").Value = 46248.332
$ws.Range("J126").Value = 46248.332
$ws.Range("L126").Value = 46248.332
$ws.Range("N126").Value = -56128.332

# Row 132 - Growth Formula Lambda
$ws.Range("H132").Value = 2275.4783
$ws.Range("I132").Value = 2008.4546
$ws.Range("J132").Value = 8150
$ws.Range("K132").Value = 6025.3638
$ws.Range("L132").Value = 24450
$ws.Range("M132").Value = -3495.3638
$ws.Range("N132").Value = -29510

# Row 137 - Magnesia Whetstone
$ws.Range("H137").Value = 2587341
$ws.Range("I137").Value = 1064869.6
$ws.Range("K137").Value = 3194608.8
$ws.Range("M137").Value = -3192058.8

# Row 138 - Cunning Craftsman's Tisane
$ws.Range("H138").Value = 274725.7
$ws.Range("I138").Value = 1699.1818
$ws.Range("J138").Value = 497191.75
$ws.Range("K138").Value = 5097.5454
$ws.Range("L138").Value = 1491575.25
$ws.Range("M138").Value = 42.45460000000003
$ws.Range("N138").Value = -1501855.25

$ws = $wb.Worksheets.Item("ARM")
# Row 32 - Steel Ingot
$ws.Range("H32").Value = 2789.48
$ws.Range("I32").Value = 2485.4624
$ws.Range("J32").Value = 6828.5713
$ws.Range("K32").Value = 2485.4624
$ws.Range("L32").Value = 6828.5713
$ws.Range("M32").Value = -2198.4624
$ws.Range("N32").Value = -7402.5713

# Row 44 - Mythril Plate
$ws.Range("H44").Value = 33974.5
$ws.Range("J44").Value = 33974.5
$ws.Range("L44").Value = 33974.5
$ws.Range("N44").Value = -34950.5

# Row 55 - Mythril Elmo
$ws.Range("H55").Value = 17208.143
$ws.Range("J55").Value = 17208.143
$ws.Range("L55").Value = 17208.143
$ws.Range("N55").Value = -17838.143

# Row 80 - Titanium Hoplon
$ws.Range("H80").Value = 25111
$ws.Range("J80").Value = 31148
$ws.Range("L80").Value = 31148
$ws.Range("N80").Value = -33144

# Row 83 - Titanium Hoplon
$ws.Range("H83").Value = 25111
$ws.Range("J83").Value = 31148
$ws.Range("L83").Value = 93444
$ws.Range("N83").Value = -103428

# Row 103 - Doman Steel Greaves of Striking
$ws.Range("H103").Value = 36448
$ws.Range("J103").Value = 36448
$ws.Range("L103").Value = 36448
$ws.Range("N103").Value = -38792

$ws = $wb.Worksheets.Item("BSM")
# Row 35 - Crowsbeak Hammer
$ws.Range("H35").Value = 31300
$ws.Range("J35").Value = 31300
$ws.Range("L35").Value = 31300
$ws.Range("N35").Value = -31920

# Row 82 - Titanium Lump Hammer
$ws.Range("H82").Value = 45763.434
$ws.Range("J82").Value = 25968.732
$ws.Range("L82").Value = 25968.732
$ws.Range("N82").Value = -26734.732

# Row 85 - Titanium Lump Hammer
$ws.Range("H85").Value = 45763.434
$ws.Range("J85").Value = 25968.732
$ws.Range("L85").Value = 25968.732
$ws.Range("N85").Value = -28620.732

# Row 122 - High Durium Tathlums
$ws.Range("H122").Value = 34246.363
$ws.Range("J122").Value = 34246.363
$ws.Range("L122").Value = 34246.363
$ws.Range("N122").Value = -44046.363

$ws = $wb.Worksheets.Item("CRP")
# Row 31 - Walnut Lumber
$ws.Range("H31").Value = 1996.2069
$ws.Range("I31").Value = 1564.44
$ws.Range("J31").Value = 4694.75
$ws.Range("K31").Value = 1564.44
$ws.Range("L31").Value = 4694.75
$ws.Range("M31").Value = -1269.44
$ws.Range("N31").Value = -5284.75

# Row 34 - Walnut Lumber
$ws.Range("H34").Value = 1996.2069
$ws.Range("I34").Value = 1564.44
$ws.Range("J34").Value = 4694.75
$ws.Range("K34").Value = 1564.44
$ws.Range("L34").Value = 4694.75
$ws.Range("M34").Value = -1362.44
$ws.Range("N34").Value = -5098.75

# Row 41 - Oak Longbow
$ws.Range("H41").Value = 18578.125
$ws.Range("J41").Value = 18578.125
$ws.Range("L41").Value = 18578.125
$ws.Range("N41").Value = -19434.125

# Row 50 - Cobalt Halberd
$ws.Range("H50").Value = 9167.5
$ws.Range("J50").Value = 9167.5
$ws.Range("L50").Value = 9167.5
$ws.Range("N50").Value = -10417.5

# Row 51 - Jade Crook
$ws.Range("H51").Value = 9549.25
$ws.Range("J51").Value = 9549.25
$ws.Range("L51").Value = 9549.25
$ws.Range("N51").Value = -11021.25

# Row 60 - Yew Longbow
$ws.Range("H60").Value = 28537.8
$ws.Range("J60").Value = 28537.8
$ws.Range("L60").Value = 28537.8
$ws.Range("N60").Value = -29559.8

# Row 61 - Jade Crook
$ws.Range("H61").Value = 9549.25
$ws.Range("J61").Value = 9549.25
$ws.Range("L61").Value = 9549.25
$ws.Range("N61").Value = -10245.25

# Row 68 - Holy Cedar Composite Bow
$ws.Range("H68").Value = 17547.5
$ws.Range("J68").Value = 17547.5
$ws.Range("L68").Value = 17547.5
$ws.Range("N68").Value = -19045.5

# Row 71 - Holy Cedar Composite Bow
$ws.Range("H71").Value = 17547.5
$ws.Range("J71").Value = 17547.5
$ws.Range("L71").Value = 52642.5
$ws.Range("N71").Value = -60130.5

# Row 134 - Ceiba Lumber
$ws.Range("H134").Value = 2770.4307
$ws.Range("I134").Value = 3079.8333
$ws.Range("J134").Value = 1896.8235
$ws.Range("K134").Value = 9239.499899999999
$ws.Range("L134").Value = 5690.470499999999
$ws.Range("M134").Value = -6704.499899999999
$ws.Range("N134").Value = -10760.4705

$ws = $wb.Worksheets.Item("CUL")
# Row 5 - Maple Syrup
$ws.Range("H5").Value = 267022.4
$ws.Range("I5").Value = 324.31708
$ws.Range("J5").Value = 1481980.4
$ws.Range("K5").Value = 972.9512399999999
$ws.Range("L5").Value = 4445941.199999999
$ws.Range("M5").Value = -860.9512399999999
$ws.Range("N5").Value = -4446165.199999999

# Row 122 - Northern Sea Salt
$ws.Range("H122").Value = 55130.715
$ws.Range("I122").Value = 361.7143
$ws.Range("K122").Value = 3255.4287
$ws.Range("M122").Value = -805.4286999999999

# Row 135 - Royal Maple Syrup
$ws.Range("H135").Value = 267022.4
$ws.Range("I135").Value = 324.31708
$ws.Range("J135").Value = 1481980.4
$ws.Range("K135").Value = 2918.85372
$ws.Range("L135").Value = 13337823.6
$ws.Range("M135").Value = -383.8537199999996
$ws.Range("N135").Value = -13342893.6

$ws = $wb.Worksheets.Item("GSM")
# Row 10 - Bone Necklace
$ws.Range("H10").Value = 4680662
$ws.Range("I10").Value = 5608794.5
$ws.Range("J10").Value = 40000
$ws.Range("K10").Value = 5608794.5
$ws.Range("L10").Value = 40000
$ws.Range("M10").Value = -5608625.5
$ws.Range("N10").Value = -40338

# Row 12 - Bone Armillae
$ws.Range("M12").ClearContents()
$ws.Range("H12").Value = 1000000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1000000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1000000
$ws.Range("N12").Value = -1000280

# Row 57 - Electrum Circlet (Amber)
$ws.Range("H57").Value = 20092
$ws.Range("J57").Value = 20092
$ws.Range("L57").Value = 20092
$ws.Range("N57").Value = -21732

# Row 123 - Ametrine Ring of Fending
$ws.Range("H123").Value = 20835.928
$ws.Range("J123").Value = 20835.928
$ws.Range("L123").Value = 20835.928
$ws.Range("N123").Value = -25735.928

$ws = $wb.Worksheets.Item("LTW")
# Row 40 - Toad Leather
$ws.Range("H40").Value = 1578.8
$ws.Range("I40").Value = 1578.8
$ws.Range("K40").Value = 1578.8
$ws.Range("M40").Value = -1442.8

# Row 61 - Raptor Leather
$ws.Range("N61").ClearContents()
$ws.Range("H61").Value = 2520.8
$ws.Range("I61").Value = 2520.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2520.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2318.8

# Row 113 - Atrociraptor Leather
$ws.Range("N113").ClearContents()
$ws.Range("H113").Value = 2520.8
$ws.Range("I113").Value = 2520.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2520.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -350.8000000000002

# Row 132 - Silver Lobo Leather
$ws.Range("H132").Value = 2287.9778
$ws.Range("I132").Value = 2168.3076
$ws.Range("J132").Value = 3065.8333
$ws.Range("K132").Value = 6504.9228
$ws.Range("L132").Value = 9197.499899999999
$ws.Range("M132").Value = -3974.9228
$ws.Range("N132").Value = -14257.4999

# Row 136 - Br'aax Leather
$ws.Range("H136").Value = 1912.8857
$ws.Range("I136").Value = 1864.1724
$ws.Range("J136").Value = 2148.3333
$ws.Range("K136").Value = 5592.5172
$ws.Range("L136").Value = 6444.999899999999
$ws.Range("M136").Value = -3042.5172
$ws.Range("N136").Value = -11544.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 54 - Woolen Tights
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0

# Row 109 - Brightlinen Turban of Crafting
$ws.Range("H109").Value = 17738.5
$ws.Range("J109").Value = 17738.5
$ws.Range("L109").Value = 17738.5
$ws.Range("N109").Value = -20512.5

# Row 132 - Snow Cotton Cloth
$ws.Range("H132").Value = 1371.2338
$ws.Range("I132").Value = 1608.1818
$ws.Range("J132").Value = 778.86365
$ws.Range("K132").Value = 4824.5454
$ws.Range("L132").Value = 2336.59095
$ws.Range("M132").Value = -2294.5454
$ws.Range("N132").Value = -7396.59095

# Row 136 - Sarcenet Cloth
$ws.Range("H136").Value = 1198.6133
$ws.Range("I136").Value = 1142.8
$ws.Range("K136").Value = 3428.4
$ws.Range("M136").Value = -878.3999999999996
